# Update countries & provincias Spain
# - Reorders a few country rows (Pakistan above Japon, Colombia above
#   Grecia/Sudafrica, Guyana above Bahamas/Guinea-Bisau/Eritrea) and
#   refreshes the COVID case figures for the affected rows, plus bumps the
#   "Datos actualizados" timestamp in the header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 22:52"

# --- Helper: write a full data row (Country, Total, New, Active, ------
#             Recovered, Critical, DeathsToday, Deaths) ----------------
function Set-Row($row, $country, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $new
    $ws.Cells.Item($row, 4).Value = $active
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# --- Rows whose stats were refreshed (no reordering) -------------------
Set-Row 4   "Estados Unidos" 423046 22711 22187 386383 9234 1635 14476
Set-Row 8   "Alemania"       111779 4116  36081 73502  4895 180  2196
Set-Row 14  "Suiza"          23280  1027  9800  12585  391  74   895
Set-Row 16  "Canada"         19195  1298  4548  14220  426  46   427
Set-Row 74  "Bosnia y Herzegovina" 804 40 79 691 4 1 34
Set-Row 100 "Malta"          299    6     16    282    4    1    1
Set-Row 138 "Barbados"       63     0     8     52     4    0    3

# --- Pakistan now reported ahead of Japon with refreshed stats ---------
Set-Row 35 "Pakistan" 4263 228 467 3735 25 4 61
Set-Row 36 "Japon"    4257 0   622 3542 80 0 93

# --- Colombia now reported ahead of Grecia / Sudafrica ------------------
Set-Row 50 "Colombia"   2054 274 123 1877 76 4 54
Set-Row 51 "Grecia"     1884 52  269 1532 84 2 83
Set-Row 52 "Sudafrica"  1845 96  95  1732 7  5 18

# --- Guyana now reported ahead of Bahamas / Guinea-Bisau / Eritrea ------
Set-Row 151 "Guyana"       37 4 8 23 4 1 6
Set-Row 152 "Bahamas"      36 3 5 25 1 0 6
Set-Row 153 "Guinea-Bisau" 33 0 0 33 0 0 0
Set-Row 154 "Eritrea"      33 2 0 33 0 0 0
